$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.097.44'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '3.739.54'
$ws.Range('E3').Value = '  -1.80%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '614.03'
$ws.Range('E5').Value = '  -0.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.52'
$ws.Range('E6').Value = '  -0.93%  '
$ws.Range('D7').Value = '3.741.01'
$ws.Range('E7').Value = '  -1.67%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('E10').Value = '  -2.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.53'
$ws.Range('E11').Value = '  +2.99%  '
$ws.Range('E12').Value = '  -2.33%  '
$ws.Range('E13').Value = '  -2.69%  '
$ws.Range('E14').Value = '  -1.93%  '
$ws.Range('D15').Value = '4.363.70'
$ws.Range('E15').Value = '  -1.82%  '
$ws.Range('D16').Value = '3.740.03'
$ws.Range('D17').Value = '69.192.06'
$ws.Range('E18').Value = '  -2.81%  '
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.31'
$ws.Range('E20').Value = '  -2.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '497.92'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.29'
$ws.Range('E22').Value = '  -3.45%  '
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.75'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.79'
$ws.Range('E26').Value = '  -3.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.70'
$ws.Range('E27').Value = '  -3.89%  '
$ws.Range('E28').Value = '  -3.29%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  +0.53%  '
$ws.Range('E31').Value = '  +3.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.95'
$ws.Range('E32').Value = '  +1.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '30.44'
$ws.Range('E33').Value = '  -4.47%  '
$ws.Range('E34').Value = '  -1.70%  '
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('E36').Value = '  -0.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.08'
$ws.Range('E37').Value = '  -2.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.346'
$ws.Range('E38').Value = '  +1.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.137'
$ws.Range('E39').Value = '  +3.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '450.95'
$ws.Range('E40').Value = '  +6.86%  '
$ws.Range('E41').Value = '  -5.29%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.01'
$ws.Range('E42').Value = '  +8.97%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '49.69'
$ws.Range('E43').Value = '  -3.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '44.78'
$ws.Range('E44').Value = '  +1.25%  '
$ws.Range('E45').Value = '  -2.58%  '
$ws.Range('D46').Value = '2.938.55'
$ws.Range('E46').Value = '  -4.23%  '
$ws.Range('E47').Value = '  -2.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '27.24'
$ws.Range('E48').Value = '  -1.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '138.40'
$ws.Range('E50').Value = '  +1.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.47'
$ws.Range('E51').Value = '  -0.51%  '
